# Weekly update to the "Ají" (chili pepper) hortaliza price sheet.
# Two new weekly price records are inserted at the top of the existing
# block (rows 450-461), pushing that block down to rows 452-463.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 450-451; everything from the old row 450
# down shifts down by two rows (old 450->452, old 451->453, ..., old
# 461->463).
$ws.Range("A450:R451").EntireRow.Insert()

# New row 450: Americana (o) record for 2023-11-09
$ws.Range("A450").Value = 9
$ws.Range("B450").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C450").Value = "Metropolitana"
$ws.Range("D450").Value2 = 45239
$ws.Range("E450").Value = 13
$ws.Range("F450").Value = 100112021
$ws.Range("G450").Value = "Ají"
$ws.Range("H450").Value = "Americana (o)"
$ws.Range("I450").Value = "Primera"
$ws.Range("J450").Value = 52
$ws.Range("K450").Value = 39000
$ws.Range("L450").Value = 41000
$ws.Range("M450").Value = 40000
$ws.Range("N450").Value = '$/caja 25 kilos'
$ws.Range("O450").Value = "Provincia de Limarí"
$ws.Range("P450").Value = 1600
$ws.Range("Q450").Value = 25
$ws.Range("R450").Value = "Hortaliza"

# New row 451: Inferno record for 2023-11-09
$ws.Range("A451").Value = 9
$ws.Range("B451").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C451").Value = "Metropolitana"
$ws.Range("D451").Value2 = 45239
$ws.Range("E451").Value = 13
$ws.Range("F451").Value = 100112021
$ws.Range("G451").Value = "Ají"
$ws.Range("H451").Value = "Inferno"
$ws.Range("I451").Value = "Primera"
$ws.Range("J451").Value = 70
$ws.Range("K451").Value = 27000
$ws.Range("L451").Value = 29000
$ws.Range("M451").Value = 28000
$ws.Range("N451").Value = '$/caja 10 kilos'
$ws.Range("O451").Value = "Región de Arica y Parinacota"
$ws.Range("P451").Value = 2800
$ws.Range("Q451").Value = 10
$ws.Range("R451").Value = "Hortaliza"
